$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 238, shifting existing 238-244 down to 242-248
$ws.Range("A238:A241").EntireRow.Insert()

# Fill in the 4 new rows (238-241) for the new week (2021-09-09, serial 44448)
$rows = @(238, 239, 240, 241)
$quality = @("Especial", "Primera", "Segunda", "Tercera")
$volumen = @(200, 500, 360, 300)
$kmin = @(12500, 10500, 8500, 5500)
$lmax = @(13000, 11000, 9000, 6000)
$mavg = @(12750, 10750, 8750, 5750)
$pkg = @(708, 597, 486, 319)

for ($i = 0; $i -lt 4; $i++) {
    $r = $rows[$i]
    $ws.Cells.Item($r, 1).Value = 2
    $ws.Cells.Item($r, 2).Value = "Comercializadora del Agro de Limarí"
    $ws.Cells.Item($r, 3).Value = "Coquimbo"
    $ws.Cells.Item($r, 4).Value = 44448
    $ws.Cells.Item($r, 5).Value = 4
    $ws.Cells.Item($r, 6).Value = 100112043
    $ws.Cells.Item($r, 7).Value = "Pepino dulce"
    $ws.Cells.Item($r, 8).Value = "Cultivar IV Región"
    $ws.Cells.Item($r, 9).Value = $quality[$i]
    $ws.Cells.Item($r, 10).Value = $volumen[$i]
    $ws.Cells.Item($r, 11).Value = $kmin[$i]
    $ws.Cells.Item($r, 12).Value = $lmax[$i]
    $ws.Cells.Item($r, 13).Value = $mavg[$i]
    $ws.Cells.Item($r, 14).Value = "$/bandeja 18 kilos"
    $ws.Cells.Item($r, 15).Value = "Provincia de Limarí"
    $ws.Cells.Item($r, 16).Value = $pkg[$i]
    $ws.Cells.Item($r, 17).Value = 18
    $ws.Cells.Item($r, 18).Value = "Hortaliza"
}
